{"js": "// The edit drops the document's closing block: the \"Under the cash flow\n// additivity principle...\" paragraph, the FV2 = PV0 x (1+r2)^2 = ... formula\n// paragraph right after it, and the final \"(25)\" equation-number paragraph.\n// Those three paragraphs are deleted in their entirety, leaving the body\n// ending at the \"(13)\" paragraph that precedes them.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the \"(13)\" equation-number paragraph that should become the new last\n// paragraph of the document.\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"(13)\") {\n    markerIndex = i;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error('Could not find the \"(13)\" paragraph marking the end of the kept content.');\n}\n\n// Delete every paragraph after the marker, from the end backward so the\n// indices of not-yet-deleted items stay valid.\nfor (let i = items.length - 1; i > markerIndex; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The edit drops the document's closing block: the \"Under the cash flow\n# additivity principle...\" paragraph, the FV2 = PV0 x (1+r2)^2 = ... formula\n# paragraph right after it, and the final \"(25)\" equation-number paragraph.\n# Those three paragraphs are deleted in their entirety, leaving the body\n# ending at the \"(13)\" paragraph that precedes them.\n$d = $word.ActiveDocument\n\n# Find the \"(13)\" equation-number paragraph that should become the new last\n# paragraph of the document.\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"(13)\") {\n        $markerIndex = $i\n    }\n}\n\nif ($markerIndex -eq -1) {\n    throw 'Could not find the \"(13)\" paragraph marking the end of the kept content.'\n}\n\n# Delete every paragraph after the marker, from the end backward so earlier\n# indices stay valid while we work.\nfor ($i = $d.Paragraphs.Count; $i -gt $markerIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
